# "membuat upload data master" - turn the Agama upload template into a
# clean master-data upload sheet:
#   - keep only the header + the two real sample rows (drop the blank
#     placeholder rows 5-12)
#   - normalise the two sample "Tgl Input" dates to the same date
#   - re-apply the (re-measured) column widths
#   - leave the saved selection covering the whole used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalise the sample dates in D3/D4 to 12 Dec 2022
$ws.Range("D3").Value = [datetime]"2022-12-12"
$ws.Range("D4").Value = [datetime]"2022-12-12"

# Drop the empty placeholder rows (old rows 5 through 12); only the
# header row and the two sample data rows remain
$null = $ws.Range("A5:F12").EntireRow.Delete()

# Re-measure the columns now that the template only holds real data
$ws.Columns("A").ColumnWidth = 4.5
$ws.Columns("B").ColumnWidth = 11.833333333333334
$ws.Columns("C").ColumnWidth = 18.5
$ws.Columns("D").ColumnWidth = 15.666666666666666
$ws.Columns("E").ColumnWidth = 13.333333333333334

# Leave the selection on the full used range instead of the old C6
$null = $ws.Range("A1:F4").Select()
